# Style the compute-instance performance comparison sheet:
#  - rename time headers to include units "(sec)"
#  - delete the blank column F (setup-time column shifts left into it)
#  - round the cost formulas to 3 decimal places
#  - tidy up column widths and restore the cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify the units in the header row.
$ws.Range("D1").Value = "Epoch training time(sec)"
$ws.Range("E1").Value = "Epoch testing time(sec)"
$ws.Range("G1").Value = "setup time(sec)"

# Column F was an empty spacer column; removing it shifts the
# "setup time" / "training cost per epoch" / "total cost" columns
# left by one (G->F, H->G, I->H) while preserving their per-cell styles.
$ws.Columns.Item(6).Delete()

# The cost formulas are now rounded to 3 decimal places.
For ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=ROUND(D" + $r + "*C" + $r + "/3600,3)"
    $ws.Cells.Item($r, 8).Formula = "=ROUND(B" + $r + "*C" + $r + "/3600,3)"
}

# Tidy column widths to fit the new headers/content.
$ws.Columns.Item(5).ColumnWidth = 21.82
$ws.Columns.Item(6).ColumnWidth = 14.16
$ws.Columns.Item(7).ColumnWidth = 19.58
$ws.Columns.Item(8).ColumnWidth = 12.78

# Restore the selection to where the author left off.
$ws.Range("F2").Select()
